# Update cryptocurrency price (column D) and 1h volume change (column E) values
# for the rows that changed, per the upstream data refresh.
# Values are forced to remain text (matching the source inlineStr cells) so
# numeric-looking strings (e.g. "518.63", "1.00") are not auto-converted by
# Excel into floating point numbers, which would lose their exact formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue 2 4 "56.875.82"
Set-TextValue 2 5 "  +4.18%  "

# Row 3
Set-TextValue 3 4 "2.345.89"
Set-TextValue 3 5 "  +2.91%  "

# Row 4
Set-TextValue 4 5 "  -0.06%  "

# Row 5
Set-TextValue 5 4 "518.63"
Set-TextValue 5 5 "  +2.63%  "

# Row 6
Set-TextValue 6 4 "134.26"
Set-TextValue 6 5 "  +4.30%  "

# Row 7
Set-TextValue 7 5 "  +0.36%  "

# Row 8
Set-TextValue 8 5 "  +1.72%  "

# Row 9
Set-TextValue 9 4 "2.342.65"
Set-TextValue 9 5 "  +1.96%  "

# Row 10
Set-TextValue 10 5 "  +6.73%  "

# Row 11
Set-TextValue 11 5 "  -0.55%  "

# Row 12
Set-TextValue 12 4 "5.21"
Set-TextValue 12 5 "  +6.15%  "

# Row 13
Set-TextValue 13 5 "  -0.04%  "

# Row 14
Set-TextValue 14 4 "2.759.90"
Set-TextValue 14 5 "  +2.73%  "

# Row 15
Set-TextValue 15 5 "  +0.89%  "

# Row 16
Set-TextValue 16 4 "56.771.23"
Set-TextValue 16 5 "  +3.85%  "

# Row 17
Set-TextValue 17 5 "  +2.62%  "

# Row 18
Set-TextValue 18 4 "2.364.93"
Set-TextValue 18 5 "  +3.27%  "

# Row 19
Set-TextValue 19 5 "  +0.25%  "

# Row 20
Set-TextValue 20 5 "  +2.09%  "

# Row 21
Set-TextValue 21 4 "321.68"
Set-TextValue 21 5 "  +4.95%  "

# Row 22
Set-TextValue 22 4 "6.53"
Set-TextValue 22 5 "  +0.53%  "

# Row 23
Set-TextValue 23 4 "1.00"
Set-TextValue 23 5 "  +0.13%  "

# Row 24
Set-TextValue 24 4 "60.71"
Set-TextValue 24 5 "  +0.62%  "

# Row 25
Set-TextValue 25 5 "  +0.39%  "

# Row 26
Set-TextValue 26 5 "  +7.39%  "

# Row 27
Set-TextValue 27 4 "7.76"
Set-TextValue 27 5 "  +3.95%  "

# Row 28
Set-TextValue 28 4 "1.24"
Set-TextValue 28 5 "  +10.70%  "

# Row 29
Set-TextValue 29 4 "170.33"
Set-TextValue 29 5 "  -0.72%  "

# Row 30
Set-TextValue 30 5 "  +5.48%  "

# Row 31
Set-TextValue 31 5 "  +3.34%  "

# Row 32
Set-TextValue 32 4 "6.18"
Set-TextValue 32 5 "  +1.63%  "

# Row 33
Set-TextValue 33 4 "18.25"
Set-TextValue 33 5 "  +1.58%  "

# Row 34
Set-TextValue 34 5 "  +0.06%  "

# Row 35
Set-TextValue 35 5 "  +0.43%  "

# Row 36
Set-TextValue 36 4 "1.24"
Set-TextValue 36 5 "  +3.43%  "

# Row 37
Set-TextValue 37 4 "0.926"
Set-TextValue 37 5 "  +1.87%  "

# Row 38
Set-TextValue 38 5 "  +4.93%  "

# Row 39
Set-TextValue 39 5 "  +8.00%  "

# Row 40
Set-TextValue 40 4 "37.76"
Set-TextValue 40 5 "  +3.24%  "

# Row 41
Set-TextValue 41 5 "  +0.87%  "

# Row 42
Set-TextValue 42 4 "3.59"
Set-TextValue 42 5 "  +5.45%  "

# Row 43
Set-TextValue 43 4 "136.88"
Set-TextValue 43 5 "  +5.55%  "

# Row 44
Set-TextValue 44 4 "276.82"
Set-TextValue 44 5 "  +10.40%  "

# Row 45
Set-TextValue 45 4 "5.09"
Set-TextValue 45 5 "  +5.36%  "

# Row 46
Set-TextValue 46 4 "0.0935"
Set-TextValue 46 5 "  +2.89%  "

# Row 47
Set-TextValue 47 5 "  +1.01%  "

# Row 48
Set-TextValue 48 4 "0.562"
Set-TextValue 48 5 "  +2.27%  "

# Row 49
Set-TextValue 49 5 "  +5.16%  "

# Row 50
Set-TextValue 50 5 "  +1.17%  "

# Row 51
Set-TextValue 51 4 "16.80"
Set-TextValue 51 5 "  +2.34%  "
